# "more rooms added in default version"
#
# Updates the existing hotel room-pricing table and appends rows for the
# remaining rooms on floors 1-4 (rooms 105/106, and all of floors 2-4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the small summary table (row 2) ---------------------------
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 6

# --- Correct the existing room rows (floor 1, rooms 102-104) ----------
$ws.Range("B5").Value = 25
$ws.Range("C5").Value = 2

$ws.Range("B6").Value = 30
$ws.Range("C6").Value = 2

$ws.Range("B7").Value = 45
$ws.Range("C7").Value = 4

# --- Append the new rooms (row 8 through row 27) -----------------------
# Seed the new rows with the same formatting as row 7 (the last existing
# data row), then overwrite the values.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C27").PasteSpecial(-4122)

$data = @(
  @(8,  105, 50, 5),
  @(9,  106, 60, 6),
  @(10, 201, 20, 2),
  @(11, 202, 25, 2),
  @(12, 203, 30, 2),
  @(13, 204, 45, 4),
  @(14, 205, 50, 5),
  @(15, 206, 60, 6),
  @(16, 301, 20, 2),
  @(17, 302, 25, 2),
  @(18, 303, 30, 2),
  @(19, 304, 45, 4),
  @(20, 305, 50, 5),
  @(21, 306, 60, 6),
  @(22, 401, 20, 2),
  @(23, 402, 25, 2),
  @(24, 403, 30, 2),
  @(25, 404, 45, 4),
  @(26, 405, 50, 5),
  @(27, 406, 60, 6)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
}

# Floors 3 (rooms 304-306) and 4 (rooms 401-406) get their "Price Per
# Night" cells typed fresh rather than copied down, so they pick up a
# plain 2-decimal number format instead of the wrap-text style used
# above.
$ws.Range("B19:B27").ClearFormats()
$ws.Range("B19:B27").NumberFormat = "0.00"
for ($r = 19; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = $data[$r - 8][2]
}

# --- Cosmetic touch-ups to match the final sheet -----------------------
# Header rows no longer need their explicit 28.5pt height.
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(3).AutoFit()

# Selection ends up on I7.
$ws.Range("I7").Select()
